$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Insert a new column before F. This shifts the old F column (Status)
#    to G, carrying its column-width definition and the "no solution yet"
#    (cross mark) style along with it.
$ws.Columns("F:F").Insert()

# 2. New "Solutions" header in F1 (inherits the bold header style from insert)
$ws.Range("F1").Value = "Solutions"

# 3. New "Understanding" header in H1 - copy header formatting from A1 first
$ws.Range("A1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("H1").Value = "Understanding"

# 4. Number of solutions for row 2 (Contain Duplicates)
$ws.Range("F2").Value = 3

# 5. Row 4 (Two Sum) status -> not done (cross mark). G3 still carries the
#    original cross-mark style inherited from the old F3, so grab that
#    formatting for G4 before G3's own status gets updated below.
$ws.Range("G3").Copy()
$ws.Range("G4").PasteSpecial(-4122)
$ws.Range("G4").Value = "$([char]0x274C)"

# 6. Row 3 (Maximum Subarray) status -> done (checkmark), matching G2's style
$ws.Range("G2").Copy()
$ws.Range("G3").PasteSpecial(-4122)
$ws.Range("G3").Value = "$([char]0x2705)"

# 7. Understanding notes (set H3 before H2 so shared strings are interned
#    in the same order as the target workbook)
$ws.Range("H3").Value = "No sol given and didn't understood provided sol"
$ws.Range("H2").Value = "Given 1 sol and understood the other 2"

# 8. Column widths for the new columns G (12) and H (14.140625).
#    ColumnWidth is expressed in "characters" and the engine rounds the
#    stored width to the nearest 1/6th of a character, so the precise
#    inverse values below land on the intended stored widths.
$ws.Range("G1").ColumnWidth = 11.166666666666666
$ws.Range("H1").ColumnWidth = 13.307291666666666

# 9. Update the selection shown when the file is reopened
[void]$ws.Range("J11").Select()
